$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Print area grows from row 7 to row 8 -----------------------------
$printArea = $wb.Names.Item("Report!Print_Area")
$printArea.RefersTo = "='Report'!`$A`$1:`$AK`$8"

# --- 2. Build row 8 from row 6's layout (same "服務" record shape) -------
#     PasteSpecial(formats) from row 6 gives every cell style 3, except
#     M/O/P/AC which come out as style 4 (the engine's format-paste drops
#     wrapText, turning row6's P/AC style 10 into style 4 - exactly what
#     row 8 needs).
$ws.Range("A6:AK6").Copy() | Out-Null
$ws.Range("A8:AK8").PasteSpecial(-4122) | Out-Null

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "服務"
$ws.Range("C8").Value = 2025070369
$ws.Range("F8").Value = 4196
$ws.Range("G8").Value = "三重蝶愛店"
$ws.Range("H8").Value = "新北市三重區"
$ws.Range("Q8").Value = "THILF04196"
$ws.Range("R8").Value = "新北一"
$ws.Range("S8").Value = "吳宗鴻"
$ws.Range("T8").Value = 1
$ws.Range("U8").Value = "已完工"
$ws.Range("V8").Value = "2025-07-02 13:08:02"
$ws.Range("W8").Value = "2025-07-02 12:00:00"
$ws.Range("X8").Value = "2025-07-02 13:07:00"
$ws.Range("Z8").Value = 1.1
$ws.Range("AB8").Value = "到場處理"
$ws.Range("AC8").Value = "PMQ3+STAR"
$ws.Range("AD8").Value = "O"
$ws.Range("AJ8").Value = "O"
$ws.Range("AK8").Value = "O"

# --- 3. Row 7: P7 / AC7 gain wrap text (style 8 -> style 9 equivalent) ---
$ws.Range("P7").WrapText = $true
$ws.Range("AC7").WrapText = $true

# --- 4. Selection moves to AC5 -------------------------------------------
$ws.Range("AC5").Select() | Out-Null
